$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Valor Mora" (column F) values for the "Periodo Mora" 2209 (row 16) and
# 2204 (row 21) swap: the 34666 figure moves from period 2209 to period 2204,
# and period 2209 now carries the standard 40000 value.
$ws.Range("F16").Value = 40000
$ws.Range("F21").Value = 34666
